$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddress, $value) {
    $rng = $ws.Range($rangeAddress)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

# Row 2 - Bitcoin
Set-TextValue "D2" "26.673.40"
Set-TextValue "E2" "  +0.26%  "

# Row 3 - Ethereum
Set-TextValue "D3" "1.595.51"
Set-TextValue "E3" "  +0.61%  "

# Row 5 - BNB
Set-TextValue "D5" "210.96"
Set-TextValue "E5" "  +0.02%  "

# Row 6 - XRP
Set-TextValue "D6" "0.512"
Set-TextValue "E6" "  +1.16%  "

# Row 7 - USDC
Set-TextValue "E7" "  +0.03%  "

# Row 8 - Dogecoin
Set-TextValue "D8" "0.0616"
Set-TextValue "E8" "  -0.08%  "

# Row 9 - Cardano
Set-TextValue "E9" "  -1.45%  "

# Row 10 - Solana
Set-TextValue "D10" "19.46"
Set-TextValue "E10" "  -0.77%  "

# Row 11 - TRON
Set-TextValue "D11" "0.0843"
Set-TextValue "E11" "  +1.02%  "

# Row 12 - WrappedliquidstakedEther2.0
Set-TextValue "D12" "1.819.25"
Set-TextValue "E12" "  +0.63%  "

# Row 13 - WrappedEther
Set-TextValue "D13" "1.601.99"
Set-TextValue "E13" "  +0.95%  "

# Row 14 - Polkadot
Set-TextValue "E14" "  -0.27%  "

# Row 15 - Polygon
Set-TextValue "D15" "0.521"
Set-TextValue "E15" "  -1.15%  "

# Row 16 - Litecoin
Set-TextValue "D16" "64.55"
Set-TextValue "E16" "  -0.06%  "

# Row 17 - WrappedBTC
Set-TextValue "D17" "26.654.25"
Set-TextValue "E17" "  +0.25%  "

# Row 18 - ShibaInu
Set-TextValue "D18" "0.0₃0729"
Set-TextValue "E18" "  +0.06%  "

# Row 19 - Dai
Set-TextValue "E19" "  +0.05%  "

# Row 20 - BitcoinCash
Set-TextValue "D20" "207.75"
Set-TextValue "E20" "  +0.35%  "

# Row 21 - Chainlink
Set-TextValue "E21" "  +0.76%  "

# Row 22 - Uniswap
Set-TextValue "D22" "4.23"
Set-TextValue "E22" "  -0.20%  "

# Row 23 - Toncoin
Set-TextValue "E23" "  -1.40%  "

# Row 24 - Avalanche
Set-TextValue "D24" "8.84"
Set-TextValue "E24" "  -0.38%  "

# Row 25 - Monero
Set-TextValue "D25" "145.43"
Set-TextValue "E25" "  -1.13%  "

# Row 26 - BinanceUSD
Set-TextValue "E26" "  +0.18%  "

# Row 27 - Cosmos
Set-TextValue "D27" "7.19"
Set-TextValue "E27" "  -2.28%  "

# Row 28 - Stellar
Set-TextValue "E28" "  +0.90%  "

# Row 29 - EthereumClassic
Set-TextValue "E29" "  -0.27%  "

# Row 30 - Hedera
Set-TextValue "E30" "  -0.01%  "

# Row 31 - PancakeSwap
Set-TextValue "D31" "1.16"
Set-TextValue "E31" "  +0.24%  "

# Row 32 - Filecoin
Set-TextValue "E32" "  -0.57%  "

# Row 33 - ImmutableX
Set-TextValue "D33" "0.659"
Set-TextValue "E33" "  -0.41%  "

# Row 35 - Maker
Set-TextValue "D35" "1.279.76"
Set-TextValue "E35" "  -3.97%  "

# Row 36 - HuobiToken
Set-TextValue "E36" "  +1.49%  "

# Row 37 - LidoDAOToken
Set-TextValue "E37" "  -1.07%  "

# Row 38 - VeChain
Set-TextValue "E38" "  -0.56%  "

# Row 39 - ARBITRUM
Set-TextValue "D39" "0.840"
Set-TextValue "E39" "  +1.69%  "

# Row 40 - PaxDollar
Set-TextValue "E40" "  +0.03%  "

# Row 41 - FraxShare
Set-TextValue "D41" "5.41"
Set-TextValue "E41" "  +1.11%  "

# Row 43 - TrustWalletToken
Set-TextValue "E43" "  +0.36%  "

# Row 44 - Aave
Set-TextValue "D44" "63.46"
Set-TextValue "E44" "  +0.01%  "

# Row 45 - RocketPoolETH
Set-TextValue "D45" "1.731.90"
Set-TextValue "E45" "  +0.65%  "

# Row 46 - WEMIXToken
Set-TextValue "E46" "  +8.70%  "

# Row 47 - Quant
Set-TextValue "D47" "90.03"
Set-TextValue "E47" "  +0.23%  "

# Row 48 - RenderToken
Set-TextValue "D48" "1.59"
Set-TextValue "E48" "  -0.87%  "

# Row 49 - Algorand
Set-TextValue "E49" "  +2.52%  "

# Row 51 - USDD -> EnergySwap
Set-TextValue "B51" "EnergySwap"
Set-TextValue "C51" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D51" "7.48"
Set-TextValue "E51" "  +0.09%  "
